# The commit swaps the two theme parts of the deck: the presentation's
# active theme (currently the "Integral" colour set, stored in
# ppt/theme/theme2.xml and used by the slide master / the whole deck)
# is replaced with the classic "Office Theme" colour set that used to
# live in ppt/theme/theme1.xml (and vice-versa for the font/format
# schemes, which are already identical between the two themes).
#
# The PowerPoint object model exposes the live theme colour scheme via
# Slide.ThemeColorScheme (it resolves to the one-and-only slide master's
# theme for this deck). Re-pointing each of the twelve theme colour
# slots to the "Office" RGB values reproduces the effective colour swap.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

function ToRGB([int]$r, [int]$g, [int]$b) {
    return $r + ($g * 256) + ($b * 65536)
}

# Slot order: 1 dk1, 2 lt1, 3 dk2, 4 lt2, 5-10 accent1-6, 11 hlink, 12 folHlink
$officeColors = @(
    (ToRGB 0x00 0x00 0x00),  # dk1
    (ToRGB 0xFF 0xFF 0xFF),  # lt1
    (ToRGB 0x44 0x54 0x6A),  # dk2
    (ToRGB 0xE7 0xE6 0xE6),  # lt2
    (ToRGB 0x5B 0x9B 0xD5),  # accent1
    (ToRGB 0xED 0x7D 0x31),  # accent2
    (ToRGB 0xA5 0xA5 0xA5),  # accent3
    (ToRGB 0xFF 0xC0 0x00),  # accent4
    (ToRGB 0x44 0x72 0xC4),  # accent5
    (ToRGB 0x70 0xAD 0x47),  # accent6
    (ToRGB 0x05 0x63 0xC1),  # hlink
    (ToRGB 0x95 0x4F 0x72)   # folHlink
)

for ($i = 1; $i -le 12; $i++) {
    $tcs.Colors($i).RGB = $officeColors[$i - 1]
}
